$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply text number format to C2:C9 first so the values are stored as text
# (preserves leading zeros, matches numFmtId 49 "@" added to cellXfs)
$ws.Range("C2:C9").NumberFormat = "@"

# Update EMSO values in column C (rows 2-9)
$ws.Range("C2").Value = "0203952500137"
$ws.Range("C3").Value = "1308959500124"
$ws.Range("C4").Value = "2203962505231"
$ws.Range("C5").Value = "1809955500218"
$ws.Range("C6").Value = "2710963500313"
$ws.Range("C7").Value = "3107964505276"
$ws.Range("C8").Value = "2811000500017"
$ws.Range("C9").Value = "1402001505453"

# Move selection
$ws.Range("I26").Select()
